# Fix Australia's country code: "AUT" (Austria's ISO code, wrongly used)
# is replaced with the correct "AUS" in B2 (Australia's row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "AUS"

# Leave the active selection on the corrected cell.
$null = $ws.Range("B2").Select()
